$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Sam"
$ws.Range("B3").Value = 45678

$ws.Range("A4").Value = "user"
$ws.Range("B4").Value = 7896141

$ws.Range("C4").Select()
